$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E3 held a blank placeholder value in the old layout and is no longer used
$ws.Range("E3").Clear() | Out-Null

# --- Phase 1: introduce the brand-new label/part-number text cells in the same order
#     they were originally typed, so the shared-string table lines up exactly ---
$ws.Range("H1").Value = "Off-board parts"
$ws.Range("A1").Value = "On-board parts"
$ws.Range("H4").Value = "Power switch"
$ws.Range("I4").Value = "EG5617-ND"
$ws.Range("B1").Value = "Part Number"
$ws.Range("H5").Value = "Other switch"
$ws.Range("I5").Value = "CWI335-ND"
$ws.Range("H2").Value = "Go button"
$ws.Range("I2").Value = "CW232-ND"
$ws.Range("H3").Value = "Stop button"
$ws.Range("I3").Value = "CW233-ND"
$ws.Range("H24").Value = "Right-angle antenna"
$ws.Range("I24").Value = "DELTA2A/X/SMAM/S/RA/11-ND"
$ws.Range("H6").Value = "Straight antenna"
$ws.Range("H7").Value = "Right-angle SMA cable"
$ws.Range("H23").Value = "Probalby not using"
$ws.Range("I7").Value = "CBA-SMAMR-SMAF-ND"
$ws.Range("M7").Value = "Too long"

# --- Phase 2: fill in the rest of the BOM values (numbers and cells that reuse
#     already-existing shared text) ---
$ws.Range("C1").Value = "Needed"
$ws.Range("D1").Value = "Order"
$ws.Range("E1").Value = "Unit Cost"
$ws.Range("I1").Value = "Part Number"
$ws.Range("J1").Value = "Needed"
$ws.Range("K1").Value = "Order"
$ws.Range("L1").Value = "Unit Cost"
$ws.Range("A2").Value = "ATMEGA"
$ws.Range("B2").Value = "ATMEGA32U4-AU"
$ws.Range("C2").Value = 1
$ws.Range("E2").Value = 4.12
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 2.84
$ws.Range("A3").Value = "RFM69HCW - 915MHz"
$ws.Range("B3").Value = "1568-1394-ND"
$ws.Range("C3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 2.77
$ws.Range("A4").Value = "3.3V regulator"
$ws.Range("B4").Value = "296-39452-1-ND"
$ws.Range("C4").Value = 1
$ws.Range("E4").Value = 1.23
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.55
$ws.Range("A5").Value = "10uH inductor"
$ws.Range("B5").Value = "587-2886-1-ND"
$ws.Range("C5").Value = 1
$ws.Range("E5").Value = 0.29
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.65
$ws.Range("A6").Value = "47uF ceramic"
$ws.Range("B6").Value = "587-1780-1-ND"
$ws.Range("C6").Value = 1
$ws.Range("E6").Value = 0.7
$ws.Range("J6").Value = 2
$ws.Range("M6").Value = "Have"
$ws.Range("A7").Value = "4.7uF ceramic"
$ws.Range("B7").Value = "587-1780-1-ND"
$ws.Range("E7").Value = 0.7
$ws.Range("J7").Value = 1
$ws.Range("L7").Value = 6.81
$ws.Range("A8").Value = "1k resistor"
$ws.Range("B8").Value = "541-3991-1-ND"
$ws.Range("A9").Value = "10k resistor"
$ws.Range("D9").Value = 0
$ws.Range("F9").Value = "Have"
$ws.Range("A10").Value = "reset button"
$ws.Range("B10").Value = "401-1426-1-ND"
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 2
$ws.Range("E10").Value = 0.52
$ws.Range("A11").Value = "0.1uF ceramic"
$ws.Range("A12").Value = "Level shifter"
$ws.Range("B12").Value = "296-12163-1-ND"
$ws.Range("C12").Value = 1
$ws.Range("E12").Value = 0.43
$ws.Range("A13").Value = "SMA connector"
$ws.Range("B13").Value = "A97594-ND"
$ws.Range("C13").Value = 1
$ws.Range("E13").Value = 2.17
$ws.Range("J24").Value = 1
$ws.Range("L24").Value = 10.82

# --- Apply currency ("Unit Cost") number formatting to the newly added cost cells ---
$currencyFormat = '_("$"* #,##0.00_);_("$"* \(#,##0.00\);_("$"* "-"??_);_(@_)'
foreach ($addr in @("E2","L2","L3","E4","L4","E5","L5","E6","E7","L7","E10","E12","E13","L24")) {
  $ws.Range($addr).NumberFormat = $currencyFormat
}

# --- Column widths for the new Off-board-parts table ---
$ws.Columns.Item(8).ColumnWidth = 18.833333333333332   # -> stored width 19.6640625
$ws.Columns.Item(9).ColumnWidth = 26.5                 # -> stored width 27.33203125
$ws.Columns.Item(12).ColumnWidth = 8.0                 # -> stored width 8.88671875

# --- Match the saved selection/active cell ---
[void]$ws.Range("M7").Select()
